$wb = $excel.ActiveWorkbook

# --- Sheet: Estadistica_general ---
$ws1 = $wb.Worksheets.Item("Estadistica_general")
$ws1.Range("B2").Value = 0.1546574809485478
$ws1.Range("B3").Value = 2.69844014698019
$ws1.Range("B4").Value = 0.3838164737673087
$ws1.Range("B5").Value = 0.1146786765760137
$ws1.Range("B6").Value = 0.1014758454536874
$ws1.Range("B7").Value = 0.1881324494120208
$ws1.Range("B8").Value = 0.08665660395833334
$ws1.Range("B9").Value = 168

# --- Sheet: Estadistica_mensual ---
$ws2 = $wb.Worksheets.Item("Estadistica_mensual")
$ws2.Range("B2").Value = 0.3568910137159959
$ws2.Range("C2").Value = 0.1453278273961506
$ws2.Range("D2").Value = 0.2220587393671027
$ws2.Range("E2").Value = 0.1644712095783876
$ws2.Range("F2").Value = 0.1374578152580881
$ws2.Range("G2").Value = 0.1080535540374005
$ws2.Range("H2").Value = 0.1215267293278492
$ws2.Range("I2").Value = 0.1473215805147113
$ws2.Range("J2").Value = 0.1797654493497215
$ws2.Range("K2").Value = 0.1584748750842526
$ws2.Range("L2").Value = 0.1560010384798874
$ws2.Range("M2").Value = 0.1971368763580774
$ws2.Range("B3").Value = -37.23346397844076
$ws2.Range("C3").Value = 1.929977275574772
$ws2.Range("D3").Value = 6.093962034303672
$ws2.Range("E3").Value = 3.412888721896989
$ws2.Range("F3").Value = 4.013397397973781
$ws2.Range("G3").Value = 2.131294830070397
$ws2.Range("H3").Value = 2.818463263969693
$ws2.Range("I3").Value = 2.047938277630508
$ws2.Range("J3").Value = 1.642755691977874
$ws2.Range("K3").Value = 4.431015219582292
$ws2.Range("L3").Value = -4.607035486040331
$ws2.Range("M3").Value = -13.99132960573672
$ws2.Range("B4").Value = -0.3606211175089777
$ws2.Range("C4").Value = 0.2450242050720917
$ws2.Range("D4").Value = 0.8906150681839287
$ws2.Range("E4").Value = 0.3565683628701244
$ws2.Range("F4").Value = 0.3603559355420568
$ws2.Range("G4").Value = 0.5757036195007328
$ws2.Range("H4").Value = 0.1508638490301752
$ws2.Range("I4").Value = 0.3434845331076297
$ws2.Range("J4").Value = 0.5243919769302634
$ws2.Range("K4").Value = 0.2907313246883527
$ws2.Range("L4").Value = -0.5794078499363275
$ws2.Range("M4").Value = -0.3267157786154967
$ws2.Range("B5").Value = 0.2477643626985334
$ws2.Range("C5").Value = 0.138452740435
$ws2.Range("D5").Value = 0.1372150052578333
$ws2.Range("E5").Value = 0.138111324214045
$ws2.Range("F5").Value = 0.09737808610366028
$ws2.Range("G5").Value = 0.07780243492082187
$ws2.Range("H5").Value = 0.1019398689425263
$ws2.Range("I5").Value = 0.1198491193542429
$ws2.Range("J5").Value = 0.1595136905382143
$ws2.Range("K5").Value = 0.1228326968775428
$ws2.Range("L5").Value = 0.12387904794792
$ws2.Range("M5").Value = 0.117443676271025
$ws2.Range("B6").Value = 0.2334691007901333
$ws2.Range("C6").Value = 0.138452740435
$ws2.Range("D6").Value = 0.1314584328867
$ws2.Range("E6").Value = 0.137052581246275
$ws2.Range("F6").Value = 0.08374253371708258
$ws2.Range("G6").Value = 0.06144377617701559
$ws2.Range("H6").Value = 0.0839969331419684
$ws2.Range("I6").Value = 0.1059734743378429
$ws2.Range("J6").Value = 0.1375630713184999
$ws2.Range("K6").Value = 0.1122631173646285
$ws2.Range("L6").Value = 0.09992071377368003
$ws2.Range("M6").Value = 0.101138123315975
$ws2.Range("B7").Value = 0.2903333471234666
$ws2.Range("C7").Value = 0.1857777866016667
$ws2.Range("D7").Value = 0.1956759352200333
$ws2.Range("E7").Value = 0.209054176596275
$ws2.Range("F7").Value = 0.1588642932599398
$ws2.Range("G7").Value = 0.1570234449582656
$ws2.Range("H7").Value = 0.1651403587209158
$ws2.Range("I7").Value = 0.2281785822664143
$ws2.Range("J7").Value = 0.2623571553185
$ws2.Range("K7").Value = 0.1911428662217714
$ws2.Range("L7").Value = 0.19735000937368
$ws2.Range("M7").Value = 0.190875009065975
$ws2.Range("B8").Value = 0.05686424633333333
$ws2.Range("C8").Value = 0.04732504616666666
$ws2.Range("D8").Value = 0.06421750233333332
$ws2.Range("E8").Value = 0.07200159534999999
$ws2.Range("F8").Value = 0.07512175954285717
$ws2.Range("G8").Value = 0.09557966878125
$ws2.Range("H8").Value = 0.08114342557894738
$ws2.Range("I8").Value = 0.1222051079285714
$ws2.Range("J8").Value = 0.124794084
$ws2.Range("K8").Value = 0.07887974885714286
$ws2.Range("L8").Value = 0.0974292956
$ws2.Range("M8").Value = 0.08973688575
$ws2.Range("B9").Value = 3
$ws2.Range("C9").Value = 6
$ws2.Range("D9").Value = 9
$ws2.Range("E9").Value = 20
$ws2.Range("F9").Value = 35
$ws2.Range("G9").Value = 32
$ws2.Range("H9").Value = 19
$ws2.Range("I9").Value = 14
$ws2.Range("J9").Value = 14
$ws2.Range("K9").Value = 7
$ws2.Range("L9").Value = 5
$ws2.Range("M9").Value = 4

# --- Sheet: Estadistica_anual ---
$ws3 = $wb.Worksheets.Item("Estadistica_anual")
$ws3.Range("M1").Copy($ws3.Range("N1"))
$ws3.Range("N1").Value = 2018
$ws3.Range("B2").Value = 0.1500225716411334
$ws3.Range("C2").Value = 0.1291989047355318
$ws3.Range("D2").Value = 0.0818209130353407
$ws3.Range("F2").Value = 0.2233715424522018
$ws3.Range("G2").Value = 0.101350728688392
$ws3.Range("H2").Value = 0.1497223002801608
$ws3.Range("I2").Value = 0.1849151099518127
$ws3.Range("J2").Value = 0.1500274869660959
$ws3.Range("K2").Value = 0.1876924827190067
$ws3.Range("L2").Value = 0.1745611112180894
$ws3.Range("M2").Value = 0.112390263447821
$ws3.Range("B3").Value = -1.238391763450412
$ws3.Range("C3").Value = 3.347511195450873
$ws3.Range("D3").Value = 1.440677442701323
$ws3.Range("F3").Value = 2.024570434822136
$ws3.Range("G3").Value = 3.742577933778201
$ws3.Range("H3").Value = 2.475415283470784
$ws3.Range("I3").Value = 5.273271626555526
$ws3.Range("J3").Value = 4.506176626286087
$ws3.Range("K3").Value = -3.872617940160398
$ws3.Range("L3").Value = 3.80782949722114
$ws3.Range("M3").Value = 2.404585807621899
$ws3.Range("B4").Value = -0.01068254766301596
$ws3.Range("C4").Value = 0.6522466681189094
$ws3.Range("D4").Value = 0.8860591713718771
$ws3.Range("F4").Value = 0.6244967826075859
$ws3.Range("G4").Value = 0.4353093830697635
$ws3.Range("H4").Value = 0.4189813231737837
$ws3.Range("I4").Value = 0.06844480151763736
$ws3.Range("J4").Value = 0.3819182184362902
$ws3.Range("K4").Value = -0.08604014901087036
$ws3.Range("L4").Value = 0.3676900803241546
$ws3.Range("M4").Value = 0.6186775554778738
$ws3.Range("B5").Value = 0.1259882946827857
$ws3.Range("C5").Value = 0.09950151858346819
$ws3.Range("D5").Value = 0.06318211090658572
$ws3.Range("F5").Value = 0.18042866569032
$ws3.Range("G5").Value = 0.07405489441831054
$ws3.Range("H5").Value = 0.1211650043772964
$ws3.Range("I5").Value = 0.1473073467146278
$ws3.Range("J5").Value = 0.1184295432666857
$ws3.Range("K5").Value = 0.1309095773087286
$ws3.Range("L5").Value = 0.1142294369287682
$ws3.Range("M5").Value = 0.09687040125795714
$ws3.Range("B6").Value = 0.1218058895460714
$ws3.Range("C6").Value = 0.09223634412035911
$ws3.Range("D6").Value = 0.05386465207512858
$ws3.Range("F6").Value = 0.16467146517328
$ws3.Range("G6").Value = 0.05545729772635262
$ws3.Range("H6").Value = 0.09209726328363158
$ws3.Range("I6").Value = 0.12926765635735
$ws3.Range("J6").Value = 0.1184295432666857
$ws3.Range("K6").Value = 0.1276981823774333
$ws3.Range("L6").Value = 0.1001222548887409
$ws3.Range("M6").Value = 0.08642464201718573
$ws3.Range("B7").Value = 0.2270476298317857
$ws3.Range("C7").Value = 0.1687310686203591
$ws3.Range("D7").Value = 0.1840476277894143
$ws3.Range("F7").Value = 0.29840001417328
$ws3.Range("G7").Value = 0.1283157955684579
$ws3.Range("H7").Value = 0.1889966756436316
$ws3.Range("I7").Value = 0.2118240841351278
$ws3.Range("J7").Value = 0.1981071522666857
$ws3.Range("K7").Value = 0.2064523907583857
$ws3.Range("L7").Value = 0.1788484933432864
$ws3.Range("M7").Value = 0.1657857221600429
$ws3.Range("B8").Value = 0.1052417402857143
$ws3.Range("C8").Value = 0.0764947245
$ws3.Range("D8").Value = 0.1301829757142857
$ws3.Range("F8").Value = 0.133728549
$ws3.Range("G8").Value = 0.07285849784210525
$ws3.Range("H8").Value = 0.09689941235999999
$ws3.Range("I8").Value = 0.08255642777777779
$ws3.Range("J8").Value = 0.07967760900000001
$ws3.Range("K8").Value = 0.07875420838095237
$ws3.Range("L8").Value = 0.07872623845454545
$ws3.Range("M8").Value = 0.07936108014285713
$ws3.Range("B9").Value = 14
$ws3.Range("C9").Value = 22
$ws3.Range("D9").Value = 7
$ws3.Range("F9").Value = 5
$ws3.Range("G9").Value = 19
$ws3.Range("H9").Value = 25
$ws3.Range("I9").Value = 18
$ws3.Range("J9").Value = 7
$ws3.Range("K9").Value = 21
$ws3.Range("L9").Value = 22
$ws3.Range("M9").Value = 7
